$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Go Kinetic by Windstream"
$ws.Range("C2").Value = 4.65
$ws.Range("D2").Value = 82788
$ws.Range("F2").Value = 61742
$ws.Range("H2").Value = 4.5
$ws.Range("I2").Value = 21046
$ws.Range("B3").Value = "My altafiber"
$ws.Range("C3").Value = 4.5
$ws.Range("D3").Value = 4743
$ws.Range("E3").Value = 4.8
$ws.Range("F3").Value = 4426
$ws.Range("H3").Value = 4.2
$ws.Range("I3").Value = 317
$ws.Range("B4").Value = "Spectrum Access: Enabled Media"
$ws.Range("C4").Value = 4.5
$ws.Range("D4").Value = 265
$ws.Range("E4").Value = 4.5
$ws.Range("F4").Value = 123
$ws.Range("H4").Value = 4.5
$ws.Range("I4").Value = 142
$ws.Range("B5").Value = "Cox App"
$ws.Range("C5").Value = 4.55
$ws.Range("D5").Value = 489344
$ws.Range("E5").Value = 4.6
$ws.Range("F5").Value = 399031
$ws.Range("G5").Value = 117
$ws.Range("H5").Value = 4.5
$ws.Range("I5").Value = 90313
$ws.Range("B6").Value = "My Verizon"
$ws.Range("C6").Value = 4.6
$ws.Range("D6").Value = 5750786
$ws.Range("E6").Value = 4.6
$ws.Range("F6").Value = 4754388
$ws.Range("G6").Value = 9
$ws.Range("H6").Value = 4.6
$ws.Range("I6").Value = 996398
$ws.Range("B7").Value = "myCricket App"
$ws.Range("C7").Value = 4.3
$ws.Range("D7").Value = 233477
$ws.Range("E7").Value = 4.5
$ws.Range("F7").Value = 43149
$ws.Range("G7").Value = 90
$ws.Range("H7").Value = 4.1
$ws.Range("I7").Value = 190328
$ws.Range("B8").Value = "MyDISH Account"
$ws.Range("C8").Value = 4.1
$ws.Range("D8").Value = 381866
$ws.Range("E8").Value = 4.6
$ws.Range("F8").Value = 322237
$ws.Range("G8").Value = $null
$ws.Range("H8").Value = 3.6
$ws.Range("I8").Value = 59629
$ws.Range("B9").Value = "T-Mobile"
$ws.Range("C9").Value = 4.2
$ws.Range("D9").Value = 3187389
$ws.Range("E9").Value = 4.8
$ws.Range("F9").Value = 2412812
$ws.Range("G9").Value = 4
$ws.Range("H9").Value = 3.6
$ws.Range("I9").Value = 774577
$ws.Range("B10").Value = "Xfinity Mobile"
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = 55397
$ws.Range("E10").Value = 4.8
$ws.Range("F10").Value = 42796
$ws.Range("G10").Value = 87
$ws.Range("H10").Value = 3.2
$ws.Range("I10").Value = 12601
$ws.Range("B11").Value = "Spectrum News: Local Stories"
$ws.Range("C11").Value = 4.4
$ws.Range("D11").Value = 30351
$ws.Range("E11").Value = 4.7
$ws.Range("F11").Value = 24958
$ws.Range("G11").Value = 51
$ws.Range("H11").Value = 4.1
$ws.Range("I11").Value = 5393
$ws.Range("B12").Value = "My Spectrum"
$ws.Range("C12").Value = 4.699999999999999
$ws.Range("D12").Value = 2645951
$ws.Range("F12").Value = 2010064
$ws.Range("G12").Value = 11
$ws.Range("H12").Value = 4.6
$ws.Range("I12").Value = 635887
$ws.Range("B13").Value = "My Sprint Mobile"
$ws.Range("C13").Value = 4.3
$ws.Range("D13").Value = 1190791
$ws.Range("E13").Value = 4.5
$ws.Range("F13").Value = 1047309
$ws.Range("H13").Value = 4.1
$ws.Range("I13").Value = 143482
$ws.Range("B14").Value = "Verizon My Fios"
$ws.Range("C14").Value = 4.4
$ws.Range("D14").Value = 334805
$ws.Range("E14").Value = 4.5
$ws.Range("F14").Value = 258060
$ws.Range("G14").Value = 111
$ws.Range("H14").Value = 4.3
$ws.Range("I14").Value = 76745
$ws.Range("B15").Value = "My CenturyLink"
$ws.Range("C15").Value = 4.25
$ws.Range("D15").Value = 171396
$ws.Range("E15").Value = 4.4
$ws.Range("F15").Value = 116043
$ws.Range("G15").Value = 196
$ws.Range("H15").Value = 4.1
$ws.Range("I15").Value = 55353
$ws.Range("B16").Value = "Visible mobile"
$ws.Range("C16").Value = 4.300000000000001
$ws.Range("D16").Value = 106650
$ws.Range("E16").Value = 4.4
$ws.Range("F16").Value = 77035
$ws.Range("G16").Value = 156
$ws.Range("I16").Value = 29615
$ws.Range("B17").Value = "SpectrumU"
$ws.Range("C17").Value = 4.4
$ws.Range("D17").Value = 427
$ws.Range("E17").Value = 4.3
$ws.Range("F17").Value = 175
$ws.Range("B18").Value = "Spectrum TV"
$ws.Range("C18").Value = 4.65
$ws.Range("D18").Value = 600878
$ws.Range("E18").Value = 4.7
$ws.Range("F18").Value = 411629
$ws.Range("G18").Value = 47
$ws.Range("H18").Value = 4.6
$ws.Range("I18").Value = 189249
$ws.Range("B19").Value = "myAT&amp;T"
$ws.Range("C19").Value = 3.5
$ws.Range("D19").Value = 4268811
$ws.Range("E19").Value = 4
$ws.Range("F19").Value = 3991303
$ws.Range("G19").Value = 14
$ws.Range("H19").Value = 3
$ws.Range("I19").Value = 277508
$ws.Range("B20").Value = "Spectrum SportsNet: Live Games"
$ws.Range("C20").Value = 3.95
$ws.Range("D20").Value = 5565
$ws.Range("E20").Value = 4.6
$ws.Range("F20").Value = 4149
$ws.Range("G20").Value = $null
$ws.Range("H20").Value = 3.3
$ws.Range("I20").Value = 1416
$ws.Range("B21").Value = "MediacomConnect"
$ws.Range("C21").Value = 3.4
$ws.Range("D21").Value = 19121
$ws.Range("E21").Value = 3.3
$ws.Range("F21").Value = 6603
$ws.Range("H21").Value = 3.5
$ws.Range("I21").Value = 12518
$ws.Range("D22").Value = 382796
$ws.Range("I22").Value = 379635
$ws.Range("D23").Value = 81559
$ws.Range("I23").Value = 36338
$ws.Range("B24").Value = "Xfinity"
$ws.Range("C24").Value = 4
$ws.Range("D24").Value = 1030498
$ws.Range("E24").Value = 4.4
$ws.Range("F24").Value = 876975
$ws.Range("G24").Value = 5
$ws.Range("H24").Value = 3.6
$ws.Range("I24").Value = 153523
$ws.Range("B25").Value = "Google Fiber"
$ws.Range("C25").Value = 3.75
$ws.Range("D25").Value = 1444
$ws.Range("E25").Value = 3.6
$ws.Range("F25").Value = 291
$ws.Range("H25").Value = 3.9
$ws.Range("I25").Value = 1153
$ws.Range("B26").Value = "My Viasat"
$ws.Range("C26").Value = 3.9
$ws.Range("D26").Value = 4574
$ws.Range("E26").Value = 4.2
$ws.Range("F26").Value = 1620
$ws.Range("H26").Value = 3.6
$ws.Range("I26").Value = 2954
$ws.Range("B27").Value = "Armstrong"
$ws.Range("C27").Value = 3.9
$ws.Range("D27").Value = 473
$ws.Range("E27").Value = 3.7
$ws.Range("F27").Value = 14
$ws.Range("H27").Value = 4.1
$ws.Range("I27").Value = 459
$ws.Range("C28").Value = $null
$ws.Range("D28").Value = 1
$ws.Range("H28").Value = $null
$ws.Range("I28").Value = 0
$ws.Range("B29").Value = "HughesNet Mobile"
$ws.Range("C29").Value = 2.3
$ws.Range("D29").Value = 1844
$ws.Range("E29").Value = 1.8
$ws.Range("F29").Value = 307
$ws.Range("H29").Value = 2.8
$ws.Range("I29").Value = 1537
$ws.Range("B30").Value = "HT My Account"
$ws.Range("C30").Value = 2.05
$ws.Range("D30").Value = 53
$ws.Range("E30").Value = 1.6
$ws.Range("F30").Value = 20
$ws.Range("I30").Value = 33
$ws.Range("B31").Value = "Midco My Account"
$ws.Range("C31").Value = 1.85
$ws.Range("D31").Value = 287
$ws.Range("E31").Value = 1.3
$ws.Range("F31").Value = 24
$ws.Range("H31").Value = 2.4
$ws.Range("I31").Value = 263
$ws.Range("B32").Value = "Optimum Support"
$ws.Range("C32").Value = 1.95
$ws.Range("D32").Value = 2601
$ws.Range("E32").Value = 1.4
$ws.Range("F32").Value = 1375
$ws.Range("H32").Value = 2.5
$ws.Range("I32").Value = 1226
$ws.Range("B33").Value = "UScellular™ – My Account"
$ws.Range("C33").Value = 4.25
$ws.Range("D33").Value = 48891
$ws.Range("E33").Value = 4.4
$ws.Range("F33").Value = 28371
$ws.Range("B34").Value = "SE Next powered by Tivo"
$ws.Range("C34").Value = 2.35
$ws.Range("D34").Value = 31
$ws.Range("E34").Value = 2.9
$ws.Range("F34").Value = 8
$ws.Range("B35").Value = "Optimum TV"
$ws.Range("C35").Value = 4.1
$ws.Range("D35").Value = 15893
$ws.Range("E35").Value = 4.6
$ws.Range("F35").Value = 13437
$ws.Range("H35").Value = 3.6
$ws.Range("I35").Value = 2456
$ws.Range("B36").Value = "Breezeline TV"
$ws.Range("C36").Value = 1.65
$ws.Range("D36").Value = 212
$ws.Range("E36").Value = 1.2
$ws.Range("F36").Value = 38
$ws.Range("G36").Value = $null
$ws.Range("H36").Value = 2.1
$ws.Range("I36").Value = 174
$ws.Range("B37").Value = "My Blue Ridge"
$ws.Range("C37").Value = 4.5
$ws.Range("D37").Value = 4180
$ws.Range("E37").Value = 4.6
$ws.Range("F37").Value = 3298
$ws.Range("B38").Value = "myBuckeye"
$ws.Range("C38").Value = 2.8
$ws.Range("D38").Value = 101
$ws.Range("E38").Value = 2.5
$ws.Range("F38").Value = 22
$ws.Range("H38").Value = 174
$ws.Range("I38").Value = 79
